# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.252.14"
$ws.Range("E2").Value = "  +2.28%  "

$ws.Range("D3").Value = "2.498.49"
$ws.Range("E3").Value = "  +1.80%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.35%  "

$ws.Range("E7").Value = "  +1.39%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.56%  "

$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.07%  "

$ws.Range("E14").Value = "  +1.36%  "

$ws.Range("D15").Value = "2.890.34"
$ws.Range("E15").Value = "  +1.72%  "

$ws.Range("D16").Value = "2.501.90"
$ws.Range("E16").Value = "  +1.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.853"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.22%  "

$ws.Range("D18").Value = "47.187.44"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.09%  "

$ws.Range("E20").Value = "  +3.14%  "

$ws.Range("D21").Value = "0.0₃0941"
$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "247.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.95%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.04"
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "

$ws.Range("E31").Value = "  +7.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0786"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.68%  "

$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("E37").Value = "  +4.11%  "

$ws.Range("E38").Value = "  +2.84%  "

$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.34%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.60%  "

$ws.Range("E44").Value = "  +2.10%  "

$ws.Range("D45").Value = "1.991.36"
$ws.Range("E45").Value = "  +0.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.51%  "

$ws.Range("E47").Value = "  -1.62%  "

$ws.Range("E48").Value = "  -3.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.66%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.46%  "
